$wb = $excel.ActiveWorkbook

# Sheet "建物" (Building) - column I (property_category) should be "building"
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"

# Sheet "汽車" (Car) - column H (property_category) should be "car"
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
